$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns hold plain-text values (prices
# and percentage strings) rather than real numbers/percentages. Excel
# auto-converts numeric-looking input into Number/Percentage values, so each
# target cell is switched to the Text number format ("@") before the literal
# string is written, which keeps the data round-tripping as text exactly
# like the original scraped values.

function Set-TextCell($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell "D2" "316.79"
Set-TextCell "E2" "3.00%"
Set-TextCell "D3" "39.48"
Set-TextCell "E3" "2.27%"
Set-TextCell "D4" "5.138"
Set-TextCell "E4" "0.84%"
Set-TextCell "D5" "0.08178"
Set-TextCell "E5" "0.77%"
Set-TextCell "D6" "1.963"
Set-TextCell "E6" "0.02%"
Set-TextCell "D7" "8.222"
Set-TextCell "E7" "3.51%"
Set-TextCell "D8" "0.9286"
Set-TextCell "E8" "-0.04%"
Set-TextCell "D9" "0.1410"
Set-TextCell "E9" "-1.58%"
Set-TextCell "D10" "0.1997"
Set-TextCell "E10" "2.18%"
Set-TextCell "D11" "0.09017"
Set-TextCell "E11" "-0.95%"
Set-TextCell "D12" "0.03500"
Set-TextCell "E12" "-0.15%"
Set-TextCell "E13" "0.00%"
Set-TextCell "D14" "0.001400"
Set-TextCell "E14" "-0.05%"
Set-TextCell "D15" "0.005860"
Set-TextCell "E15" "-3.78%"
Set-TextCell "D16" "3.651"
Set-TextCell "E16" "-2.02%"
Set-TextCell "D17" "4.241"
Set-TextCell "E17" "1.38%"
Set-TextCell "D18" "3.162"
Set-TextCell "E18" "-8.24%"
Set-TextCell "D19" "0.3465"
Set-TextCell "E19" "0.07%"
Set-TextCell "D20" "0.1304"
Set-TextCell "E20" "0.89%"
Set-TextCell "D21" "4.849"
Set-TextCell "E21" "1.11%"
Set-TextCell "E22" "-1.00%"
Set-TextCell "D23" "0.04380"
Set-TextCell "E23" "0.32%"
Set-TextCell "D24" "0.001221"
Set-TextCell "E24" "-0.07%"
Set-TextCell "D25" "0.004787"
Set-TextCell "E25" "-0.99%"
Set-TextCell "E26" "-0.16%"
Set-TextCell "D27" "0.0004000"
Set-TextCell "E27" "-10.07%"
Set-TextCell "D39" "0.02202"
Set-TextCell "E39" "5.32%"
Set-TextCell "D40" "0.05180"
Set-TextCell "E40" "1.22%"
Set-TextCell "D41" "0.007586"
Set-TextCell "E41" "1.68%"
Set-TextCell "D42" "0.009748"
Set-TextCell "E42" "-3.72%"
Set-TextCell "D43" "0.1376"
Set-TextCell "E43" "1.31%"
Set-TextCell "E44" "-0.16%"
Set-TextCell "D45" "0.009128"
Set-TextCell "E45" "-1.50%"
Set-TextCell "D46" "0.00006403"
Set-TextCell "E46" "2.22%"
Set-TextCell "D47" "0.00000000750"
Set-TextCell "E47" "-0.14%"
Set-TextCell "D48" "0.002765"
Set-TextCell "E48" "-8.72%"
Set-TextCell "D49" "0.001200"
Set-TextCell "E49" "-25.00%"
Set-TextCell "D50" "0.00002100"
Set-TextCell "E50" "-0.14%"
Set-TextCell "D51" "0.0002000"
Set-TextCell "E51" "-0.14%"
